$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new numeric-looking text would otherwise be auto-coerced
# to a Number by Excel (losing formatting such as trailing zeros /
# leading zeros in small decimals). Force them to Text first so the
# literal string is preserved exactly as in the source data.
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"

$ws.Range("D2").Value = "28.393.37"
$ws.Range("E2").Value = "  +0.11%  "

$ws.Range("D3").Value = "1.867.35"
$ws.Range("E3").Value = "  -0.15%  "

$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "331.27"
$ws.Range("E5").Value = "  -2.27%  "

$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").Value = "0.4625"
$ws.Range("E7").Value = "  -1.67%  "

$ws.Range("D8").Value = "0.4018"
$ws.Range("E8").Value = "  +1.79%  "

$ws.Range("D9").Value = "47.94"
$ws.Range("E9").Value = "  +1.36%  "

$ws.Range("D10").Value = "0.07873"
$ws.Range("E10").Value = "  -1.59%  "

$ws.Range("D11").Value = "0.9860"
$ws.Range("E11").Value = "  -2.15%  "

$ws.Range("D12").Value = "21.31"
$ws.Range("E12").Value = "  -2.91%  "

$ws.Range("D13").Value = "1.862.86"
$ws.Range("E13").Value = "  -0.32%  "

$ws.Range("D14").Value = "5.855"
$ws.Range("E14").Value = "  -2.56%  "

$ws.Range("D15").Value = "7.005"
$ws.Range("E15").Value = "  -3.83%  "

$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.07%  "

$ws.Range("D17").Value = "88.21"
$ws.Range("E17").Value = "  -3.37%  "

$ws.Range("D18").Value = "0.06535"
$ws.Range("E18").Value = "  -0.91%  "

$ws.Range("D19").Value = "0.00001020"
$ws.Range("E19").Value = "  -2.14%  "

$ws.Range("D20").Value = "17.22"
$ws.Range("E20").Value = "  -2.78%  "

$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("D22").Value = "28.370.44"
$ws.Range("E22").Value = "  +0.03%  "

$ws.Range("D23").Value = "5.351"
$ws.Range("E23").Value = "  -2.03%  "

$ws.Range("D24").Value = "10.90"
$ws.Range("E24").Value = "  -1.52%  "

$ws.Range("D25").Value = "2.231"
$ws.Range("E25").Value = "  -2.47%  "

$ws.Range("D26").Value = "2.087.73"
$ws.Range("E26").Value = "  -0.33%  "

$ws.Range("D27").Value = "156.78"
$ws.Range("E27").Value = "  -1.96%  "

$ws.Range("D28").Value = "19.38"
$ws.Range("E28").Value = "  -2.48%  "

$ws.Range("D29").Value = "2.062"
$ws.Range("E29").Value = "  -4.03%  "

$ws.Range("D30").Value = "5.316"
$ws.Range("E30").Value = "  -3.32%  "

$ws.Range("D31").Value = "117.70"
$ws.Range("E31").Value = "  -2.24%  "

$ws.Range("D32").Value = "0.9609"
$ws.Range("E32").Value = "  -1.50%  "

$ws.Range("D33").Value = "0.09361"
$ws.Range("E33").Value = "  -1.63%  "

$ws.Range("D34").Value = "3.584"
$ws.Range("E34").Value = "  -0.23%  "

$ws.Range("D35").Value = "1.386"
$ws.Range("E35").Value = "  +0.27%  "

$ws.Range("D36").Value = "5.254"
$ws.Range("E36").Value = "  -2.01%  "

$ws.Range("D37").Value = "0.06039"
$ws.Range("E37").Value = "  -1.02%  "

$ws.Range("D38").Value = "0.02204"
$ws.Range("E38").Value = "  -3.20%  "

$ws.Range("D39").Value = "8.272"
$ws.Range("E39").Value = "  -2.49%  "

$ws.Range("D40").Value = "1.162"
$ws.Range("E40").Value = "  -1.54%  "

$ws.Range("E41").Value = "  +0.06%  "

$ws.Range("D42").Value = "0.5755"
$ws.Range("E42").Value = "  -3.81%  "

$ws.Range("D43").Value = "0.1809"
$ws.Range("E43").Value = "  -3.99%  "

$ws.Range("E44").Value = "  -3.20%  "

$ws.Range("D45").Value = "1.271"
$ws.Range("E45").Value = "  -1.38%  "

$ws.Range("D46").Value = "2.302"
$ws.Range("E46").Value = "  +13.82%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "11.93"
$ws.Range("E47").Value = "  -1.64%  "

$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").Value = "0.5435"
$ws.Range("E48").Value = "  -3.44%  "

$ws.Range("D49").Value = "0.07139"
$ws.Range("E49").Value = "  +3.11%  "

$ws.Range("D50").Value = "1.891"
$ws.Range("E50").Value = "  -3.84%  "

$ws.Range("D51").Value = "111.43"
$ws.Range("E51").Value = "  +0.05%  "
